# Applies the "eliminata la relazione SCRITTA" restructuring:
#   - ETICHETTA moves before DISCO (kept as-is, loses its own pPr/underline para format)
#   - DISCO moves right after ETICHETTA; its field order changes
#     ( ... ID_etichetta, ID_genere, ID_collezione ) instead of
#     ( ... ID_etichetta, ID_collezione, ID_genere )
#   - TRACCIA moves after DISCO (unchanged text, runs merged)
#   - AUTORE moves after TRACCIA; fields become
#     (ID, nome, cognome, IPI, ID_traccia) instead of (ID, IPI, nome, carriera, cognome)
#   - DOPPIONE field order changes to
#     (ID, progressivo, quantita, formato, condizione, ID_disco, ID_collezionista)
#   - SCRITTA paragraph is removed entirely
#   - IMMAGINE loses the duplicate " tipo, " fragment (the field list already
#     contains "tipo" later, in "percorso, tipo, ")

$d = $word.ActiveDocument

# --- 1. Replace the AUTORE..SCRITTA block (paragraphs 4-9) -----------------
# Locate the block by its start (AUTORE paragraph) and end (SCRITTA paragraph)
# so the five resulting paragraphs (ETICHETTA, DISCO, TRACCIA, AUTORE,
# DOPPIONE) land exactly where the old six (AUTORE, TRACCIA, ETICHETTA,
# DISCO, DOPPIONE, SCRITTA) used to be.

$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "AUTORE *" -and $startPara -eq $null) {
        $startPara = $d.Paragraphs.Item($i)
    }
    if ($txt -like "SCRITTA (*") {
        $endPara = $d.Paragraphs.Item($i)
    }
}

$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

$blockXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r><w:t>ETICHETTA (</w:t></w:r>
  <w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>ID</w:t></w:r>
  <w:r><w:t xml:space="preserve">, nome) </w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">DISCO </w:t></w:r>
  <w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>(ID,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> titolo_disco, anno_uscita, barcode, durata_totale, ID_etichetta, ID_genere</w:t></w:r>
  <w:r><w:t>,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>ID_collezione)</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>
  <w:r><w:t xml:space="preserve">TRACCIA </w:t></w:r>
  <w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>(ID,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> titolo, durata, ID_disco)</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>
  <w:r><w:t xml:space="preserve">AUTORE </w:t></w:r>
  <w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>(ID,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> nome,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>cognome</w:t></w:r>
  <w:r><w:t xml:space="preserve">, </w:t></w:r>
  <w:r><w:t>IPI,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> ID_traccia</w:t></w:r>
  <w:r><w:t>)</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>DOPPIONE (</w:t></w:r>
  <w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>ID</w:t></w:r>
  <w:r><w:t xml:space="preserve">, </w:t></w:r>
  <w:r><w:t>progressivo,</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve">quantita, </w:t></w:r>
  <w:r><w:t>formato, condizione</w:t></w:r>
  <w:r><w:t xml:space="preserve">, </w:t></w:r>
  <w:r><w:t>ID_disco, ID_collezionista</w:t></w:r>
  <w:r><w:t>)</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$blockRange.InsertXML($blockXml)

# --- 2. Fix IMMAGINE: drop the duplicated " tipo, " fragment ---------------
# "IMMAGINE (ID, tipo, percorso, tipo, ID_disco)" -> "IMMAGINE (ID, percorso, tipo, ID_disco)"
# Scope the Find to just the IMMAGINE paragraph and replace only the FIRST
# occurrence of " tipo, " (wdReplaceOne) so the later, still-wanted
# "percorso, tipo, " text is left untouched.
$immaginePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "IMMAGINE *") {
        $immaginePara = $d.Paragraphs.Item($i)
        break
    }
}
$immRange = $d.Range($immaginePara.Range.Start, $immaginePara.Range.End)
[void]$immRange.Find.Execute(" tipo, ", $true, $false, $false, $false, $false, $true, 1, $false, " ", 1)

Write-Output "done"
